$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 30
$ws.Range("D3").Value = 37
$ws.Range("D4").Value = 42
$ws.Range("D5").Value = 35.5
$ws.Range("D6").Value = 43
$ws.Range("D7").Value = 38
$ws.Range("D9").Value = 32
$ws.Range("D10").Value = 20
$ws.Range("D11").Value = "-"
$ws.Range("D12").Value = 18
$ws.Range("D13").Value = 33.5
$ws.Range("D15").Value = 7

$ws.Range("D16").Select()
